$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing row 75 ("01-01-2021") with revised figures
$ws.Range("B75").Value = 53484
$ws.Range("C75").Value = 1672
$ws.Range("D75").Value = 485
$ws.Range("E75").Value = 1188
$ws.Range("F75").Value = 29397
$ws.Range("G75").Value = 106
$ws.Range("H75").Value = 29290
$ws.Range("I75").Value = 10127
$ws.Range("J75").Value = 10127
$ws.Range("K75").Value = 8455
$ws.Range("L75").Value = 858
$ws.Range("M75").Value = 207
$ws.Range("N75").Value = 7390
$ws.Range("O75").Value = 3832
$ws.Range("P75").Value = -705
$ws.Range("Q75").Value = 54189
$ws.Range("R75").Value = 615
$ws.Range("S75").Value = 415
$ws.Range("T75").Value = 201
$ws.Range("U75").Value = 6211
$ws.Range("V75").Value = 6211
$ws.Range("W75").Value = 43916
$ws.Range("X75").Value = 42384
$ws.Range("Y75").Value = 1532
$ws.Range("Z75").Value = 3446

# Append new row 76 ("01-04-2021")
# Format as text first so Excel stores the literal string instead of
# auto-converting it to a date serial, then restore default formatting
# so the cell carries no explicit style (matches the rest of column A).
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "01-04-2021"
$ws.Range("A76").Style = "Normal"
$ws.Range("B76").Value = 51040
$ws.Range("C76").Value = 1582
$ws.Range("D76").Value = 454
$ws.Range("E76").Value = 1128
$ws.Range("F76").Value = 26757
$ws.Range("G76").Value = 71
$ws.Range("H76").Value = 26686
$ws.Range("I76").Value = 10284
$ws.Range("J76").Value = 10284
$ws.Range("K76").Value = 8409
$ws.Range("L76").Value = 836
$ws.Range("M76").Value = 183
$ws.Range("N76").Value = 7390
$ws.Range("O76").Value = 4008
$ws.Range("P76").Value = -2900
$ws.Range("Q76").Value = 53941
$ws.Range("R76").Value = 740
$ws.Range("S76").Value = 508
$ws.Range("T76").Value = 232
$ws.Range("U76").Value = 6225
$ws.Range("V76").Value = 6225
$ws.Range("W76").Value = 43693
$ws.Range("X76").Value = 42184
$ws.Range("Y76").Value = 1509
$ws.Range("Z76").Value = 3282
